$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6513.75
$ws.Range("J18").Value = 18333
$ws.Range("L18").Value = 18333
$ws.Range("N18").Value = -18901

$ws.Range("H28").Value = 890.3333
$ws.Range("I28").Value = 898.2222
$ws.Range("J28").Value = 878.5
$ws.Range("K28").Value = 898.2222
$ws.Range("L28").Value = 878.5
$ws.Range("M28").Value = -413.2222
$ws.Range("N28").Value = -1848.5

$ws.Range("H74").Value = 4995.2666
$ws.Range("J74").Value = 3071.4285
$ws.Range("L74").Value = 3071.4285
$ws.Range("N74").Value = -4943.4285

$ws.Range("H77").Value = 4995.2666
$ws.Range("J77").Value = 3071.4285
$ws.Range("L77").Value = 15357.1425
$ws.Range("N77").Value = -24717.1425

$ws.Range("H80").Value = 4495.7666
$ws.Range("J80").Value = 7826.9375
$ws.Range("L80").Value = 23480.8125
$ws.Range("N80").Value = -25476.8125

$ws.Range("H83").Value = 4495.7666
$ws.Range("J83").Value = 7826.9375
$ws.Range("L83").Value = 70442.4375
$ws.Range("N83").Value = -80426.4375

$ws.Range("H98").Value = 3401.7368
$ws.Range("I98").Value = 1180.4517
$ws.Range("J98").Value = 13238.857
$ws.Range("K98").Value = 1180.4517
$ws.Range("L98").Value = 13238.857
$ws.Range("M98").Value = 317.5482999999999
$ws.Range("N98").Value = -16234.857

$ws.Range("H122").Value = 3401.7368
$ws.Range("I122").Value = 1180.4517
$ws.Range("J122").Value = 13238.857
$ws.Range("K122").Value = 3541.3551
$ws.Range("L122").Value = 39716.571
$ws.Range("M122").Value = -1091.3551
$ws.Range("N122").Value = -44616.571

$ws.Range("H137").Value = 57895
$ws.Range("I137").Value = 45677
$ws.Range("K137").Value = 137031
$ws.Range("M137").Value = -134481

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21905.84
$ws.Range("I32").Value = 22526.938
$ws.Range("K32").Value = 22526.938
$ws.Range("M32").Value = -22239.938

$ws.Range("H45").Value = 2384.8928
$ws.Range("I45").Value = 1710.9
$ws.Range("J45").Value = 4069.875
$ws.Range("K45").Value = 1710.9
$ws.Range("L45").Value = 4069.875
$ws.Range("M45").Value = -1333.9
$ws.Range("N45").Value = -4823.875

$ws.Range("H61").Value = 4869
$ws.Range("I61").Value = 884.4783
$ws.Range("K61").Value = 884.4783
$ws.Range("M61").Value = -672.4783

$ws.Range("H74").Value = 408428.06
$ws.Range("I74").Value = 750786.75
$ws.Range("J74").Value = 17161
$ws.Range("K74").Value = 750786.75
$ws.Range("L74").Value = 17161
$ws.Range("M74").Value = -749912.75
$ws.Range("N74").Value = -18909

$ws.Range("H77").Value = 408428.06
$ws.Range("I77").Value = 750786.75
$ws.Range("J77").Value = 17161
$ws.Range("K77").Value = 3753933.75
$ws.Range("L77").Value = 85805
$ws.Range("M77").Value = -3749565.75
$ws.Range("N77").Value = -94541

$ws.Range("H122").Value = 2774.5557
$ws.Range("I122").Value = 2621.375
$ws.Range("K122").Value = 7864.125
$ws.Range("M122").Value = -5414.125

$ws.Range("H136").Value = 4869
$ws.Range("I136").Value = 884.4783
$ws.Range("K136").Value = 2653.4349
$ws.Range("M136").Value = -103.4349000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 848.9091
$ws.Range("I22").Value = 853.8
$ws.Range("K22").Value = 853.8
$ws.Range("M22").Value = -680.8

$ws.Range("H99").Value = 1427
$ws.Range("J99").Value = 1408
$ws.Range("L99").Value = 1408

$ws.Range("H105").Value = 2129.3
$ws.Range("I105").Value = 1892
$ws.Range("J105").Value = 3474
$ws.Range("K105").Value = 1892
$ws.Range("L105").Value = 3474
$ws.Range("M105").Value = -145
$ws.Range("N105").Value = -6968

$ws.Range("H134").Value = 3174.1333
$ws.Range("I134").Value = 2938.3076
$ws.Range("K134").Value = 8814.9228
$ws.Range("M134").Value = -6279.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1863.5264
$ws.Range("I94").Value = 1369.1111
$ws.Range("K94").Value = 1369.1111
$ws.Range("M94").Value = -918.1111000000001

$ws.Range("H107").Value = 449.6
$ws.Range("I107").Value = 311.22726
$ws.Range("J107").Value = 830.125
$ws.Range("K107").Value = 311.22726
$ws.Range("L107").Value = 830.125
$ws.Range("M107").Value = 1608.77274
$ws.Range("N107").Value = -4670.125

$ws.Range("H132").Value = 101208.2
$ws.Range("I132").Value = 101208.2
$ws.Range("K132").Value = 303624.6
$ws.Range("M132").Value = -301094.6

$ws.Range("H134").Value = 2380.2
$ws.Range("J134").Value = 5499
$ws.Range("L134").Value = 16497
$ws.Range("N134").Value = -21567

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 155.83333
$ws.Range("I15").Value = 88
$ws.Range("K15").Value = 264
$ws.Range("M15").Value = -124

$ws.Range("H21").Value = 2000
$ws.Range("J21").Value = 2000
$ws.Range("L21").Value = 6000

$ws.Range("H26").Value = 40
$ws.Range("I26").Value = 40
$ws.Range("K26").Value = 120
$ws.Range("M26").Value = 168

$ws.Range("H32").Value = 7956.478
$ws.Range("I32").Value = 7638.8335
$ws.Range("J32").Value = 9100
$ws.Range("K32").Value = 22916.5005
$ws.Range("L32").Value = 27300
$ws.Range("M32").Value = -22633.5005
$ws.Range("N32").Value = -27866

$ws.Range("H33").Value = 609.7143
$ws.Range("I33").Value = 393.6
$ws.Range("J33").Value = 1150
$ws.Range("K33").Value = 2361.6
$ws.Range("L33").Value = 6900
$ws.Range("M33").Value = -2078.6
$ws.Range("N33").Value = -7466

$ws.Range("H37").Value = 42089.176
$ws.Range("J37").Value = 42089.176
$ws.Range("L37").Value = 126267.528
$ws.Range("N37").Value = -126491.528

$ws.Range("H50").Value = 8967.5
$ws.Range("I50").Value = 25349.5
$ws.Range("J50").Value = 776.5
$ws.Range("K50").Value = 76048.5
$ws.Range("L50").Value = 2329.5
$ws.Range("M50").Value = -75567.5
$ws.Range("N50").Value = -3291.5

$ws.Range("H53").Value = 8967.5
$ws.Range("I53").Value = 25349.5
$ws.Range("J53").Value = 776.5
$ws.Range("K53").Value = 76048.5
$ws.Range("L53").Value = 2329.5
$ws.Range("M53").Value = -75567.5
$ws.Range("N53").Value = -3291.5

$ws.Range("H57").Value = 9999.883
$ws.Range("J57").Value = 9999.883
$ws.Range("L57").Value = 29999.649
$ws.Range("N57").Value = -31117.649

$ws.Range("H58").Value = 4993.5
$ws.Range("I58").Value = 4993.5
$ws.Range("K58").Value = 14980.5
$ws.Range("M58").Value = -14852.5

$ws.Range("H116").Value = 8360.571
$ws.Range("I116").Value = 2349.3333
$ws.Range("K116").Value = 7047.999899999999
$ws.Range("M116").Value = -3605.999899999999

$ws.Range("H117").Value = 601.1667
$ws.Range("I117").Value = 351.75
$ws.Range("J117").Value = 1100
$ws.Range("K117").Value = 1055.25
$ws.Range("L117").Value = 3300
$ws.Range("M117").Value = 2386.75
$ws.Range("N117").Value = -10184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8080.2354
$ws.Range("I80").Value = 1417
$ws.Range("K80").Value = 1417
$ws.Range("M80").Value = -419

$ws.Range("H83").Value = 8080.2354
$ws.Range("I83").Value = 1417
$ws.Range("K83").Value = 7085
$ws.Range("M83").Value = -2093

$ws.Range("H122").Value = 4267.3335
$ws.Range("I122").Value = 4151.4287
$ws.Range("J122").Value = 4499.143
$ws.Range("K122").Value = 12454.2861
$ws.Range("L122").Value = 13497.429
$ws.Range("M122").Value = -10004.2861
$ws.Range("N122").Value = -18397.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 200
$ws.Range("I40").Value = 200
$ws.Range("K40").Value = 200
$ws.Range("M40").Value = -64

$ws.Range("H81").Value = 122249.75
$ws.Range("J81").Value = 122249.75
$ws.Range("L81").Value = 122249.75
$ws.Range("N81").Value = -124245.75

$ws.Range("H84").Value = 122249.75
$ws.Range("J84").Value = 122249.75
$ws.Range("L84").Value = 366749.25
$ws.Range("N84").Value = -376733.25

$ws.Range("H122").Value = 2768.6924
$ws.Range("I122").Value = 2799.4783
$ws.Range("K122").Value = 8398.4349
$ws.Range("M122").Value = -5948.4349

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9163.909
$ws.Range("I81").Value = 12443.429
$ws.Range("K81").Value = 24886.858
$ws.Range("M81").Value = -23825.858

$ws.Range("H84").Value = 9163.909
$ws.Range("I84").Value = 12443.429
$ws.Range("K84").Value = 124434.29
$ws.Range("M84").Value = -119130.29

$ws.Range("H107").Value = 906.8
$ws.Range("I107").Value = 887.55554
$ws.Range("J107").Value = 935.6667
$ws.Range("K107").Value = 2662.66662
$ws.Range("L107").Value = 2807.0001
$ws.Range("M107").Value = -742.66662
$ws.Range("N107").Value = -6647.0001

$ws.Range("H113").Value = 1334.95
$ws.Range("I113").Value = 1304.3334
$ws.Range("K113").Value = 3913.0002
$ws.Range("M113").Value = -1743.0002
